# Autocorrect OCR mistakes - updated annotations, caches and user guide
$wb = $excel.ActiveWorkbook

$wsSimple = $wb.Worksheets.Item("Simple Fields")
$wsSimpleFmt = $wb.Worksheets.Item("Simple Fields - Formatted")
$wsItems = $wb.Worksheets.Item("Items")
$wsItemsFmt = $wb.Worksheets.Item("Items - Formatted")

# Vendor address: add missing space before the zip code
$wsSimple.Range("B2").Value = "1580 NW Gilman Blvd Suite 1 Issaquah WA, 98027"
$wsSimpleFmt.Range("B2").Value = "1580 NW Gilman Blvd Suite 1 Issaquah WA, 98027"

# Currency: USD -> SGD
$wsSimple.Range("I2").Value = "SGD"
$wsSimpleFmt.Range("I2").Value = "SGD"

# Item description OCR corrections
$wsItems.Range("A2").Value = "green onion Pancakes ÂY/MAf (1)"
$wsItemsFmt.Range("A2").Value = "green onion Pancakes ÂY/MAf (1)"

$wsItems.Range("A3").Value = "Pan Fried Leek Dumplings #7 (2)"
$wsItemsFmt.Range("A3").Value = "Pan Fried Leek Dumplings #7 (2)"

$wsItems.Range("A4").Value = "Pork Xiao Long Bao(10) AP])\`$E(10)"
$wsItemsFmt.Range("A4").Value = "Pork Xiao Long Bao(10) AP])\`$E(10)"

$wsItems.Range("A5").Value = "Q-BA( (5) ĦEH'L (5)"
$wsItemsFmt.Range("A5").Value = "Q-BA( (5) ĦEH'L (5)"

$wsItems.Range("A6").Value = "Chicken potstickers HÈP]`$9I5(6)"
$wsItemsFmt.Range("A6").Value = "Chicken potstickers HÈP]`$9I5(6)"

$wsItems.Range("A7").Value = "Tomato Mushroom Steamed dumpli PEÅINABEEMKK (6)"
$wsItemsFmt.Range("A7").Value = "Tomato Mushroom Steamed dumpli PEÅINABEEMKK (6)"

$wsItems.Range("A8").Value = "Zucchini shrimp dumplings ĦJU]K"
$wsItemsFmt.Range("A8").Value = "Zucchini shrimp dumplings ĦJU]K"

$wsItems.Range("A9").Value = 'beef stew nodle soup (Non Spicy "H751PJ(74k)'
$wsItemsFmt.Range("A9").Value = 'beef stew nodle soup (Non Spicy "H751PJ(74k)'

$wsItems.Range("A10").Value = "dandan noodle"
$wsItemsFmt.Range("A10").Value = "dandan noodle"

$wsItems.Range("A11").Value = "banana naan bread ¥"
$wsItemsFmt.Range("A11").Value = "banana naan bread ¥"

$wsItems.Range("A12").Value = "house made plum juice ĚUNNT"
$wsItemsFmt.Range("A12").Value = "house made plum juice ĚUNNT"
